$d = $word.ActiveDocument

# --- Change 1: paragraph 2 ("<tab>Lean Red Meat\n") -----------------------
# Split the tab character into its own (unformatted) run containing 16
# spaces, keeping "Lean Red Meat\n" in a separate run with the original
# run formatting.
$p1 = $d.Paragraphs(2)
$p1Full = $p1.Range
$tabRange = $d.Range($p1Full.Start, $p1Full.Start + 1)
$tabRange.Text = ""
$insertPoint = $d.Range($p1Full.Start, $p1Full.Start)
$spacesXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">                </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($spacesXml)

# --- Change 2: merge the "Fatty meat ..." runs into a single run ----------
$text2 = "Fatty meat or heavily marbled meat will not dry adequately, and any fat left on the meat may go rancid in storage."
$d.Content.Find.Execute($text2, $true, $false, $false, $false, $false, $true, 1, $false, $text2, 2) | Out-Null

# --- Change 3: merge the "Meat should be very cold ..." runs --------------
$text3 = "Meat should be very cold when sliced to make slicing easier.  Cut into strips 1/4-inch-thick and 1 1/2 inches wide by 4 to 8 inches long."
$d.Content.Find.Execute($text3, $true, $false, $false, $false, $false, $true, 1, $false, $text3, 2) | Out-Null

Write-Host "edits applied"
